$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the wording of the existing activity log entry in G16 (shared string
# text edited in place: "Arithmetic Unit" -> "ArithUnit.vhd" for the first
# sentence).
$ws.Range("G16").Value = "Started work on ArithUnit.vhd. Stopped to look back at notes to fully understand Arithmetic Unit before continueing"

# Row 17: new activity log entry
$ws.Range("B17").Value = 6977
$ws.Range("C17").Value = 43924
$ws.Range("D17").Value = 0.084027777777777771
$ws.Range("E17").Value = 0.11597222222222221
$ws.Range("G17").Value = "Fixing LogicUnit.vhd, Adder.vhd and ArithUnit.vhd to make it compile Quartus. Not done"

# Row 18: new activity log entry
$ws.Range("B18").Value = 6977
$ws.Range("C18").Value = 43924
$ws.Range("D18").Value = 0.11597222222222221
$ws.Range("E18").Value = 0.14166666666666666
$ws.Range("G18").Value = "Added .gitignore to ignore some of the Quartus generated files"

# Row 19: new activity log entry
$ws.Range("B19").Value = 6977
$ws.Range("C19").Value = 43924
$ws.Range("D19").Value = 0.14166666666666666
$ws.Range("E19").Value = 0.14791666666666667
$ws.Range("G19").Value = "Added .gitignore to ignore Office temporary files. Makes it annoying to accidentally include temporary files to commit"

# Move / restore the active selection to B20 (matches the author's final
# cursor position when they saved the workbook).
$ws.Range("B20").Select() | Out-Null
